# Add the two new worksheets (AdminUserPage, ManageNewsPage) after the existing LoginPage sheet.
$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item(1)

$adminSheet = $wb.Worksheets.Add($null, $loginSheet)
$adminSheet.Name = "AdminUserPage"

$newsSheet = $wb.Worksheets.Add($null, $adminSheet)
$newsSheet.Name = "ManageNewsPage"

# --- AdminUserPage content ---
$adminSheet.Range("A1").Value = "testadminuser"
$adminSheet.Range("B1").Value = "testadminpassword"

$adminFont = $adminSheet.Range("A1:B1").Font
$adminFont.Name = "Consolas"
$adminFont.Size = 10
$adminFont.Family = 3
$adminFont.Color = 16711722

$adminSheet.Columns.Item(1).ColumnWidth = 25.21875
$adminSheet.Columns.Item(2).ColumnWidth = 21.77734375

# --- ManageNewsPage content ---
$newsSheet.Range("A1").Value = "Political test issues"

$newsFont = $newsSheet.Range("A1").Font
$newsFont.Name = "Consolas"
$newsFont.Size = 10
$newsFont.Family = 3
$newsFont.Color = 16711722

$newsSheet.Columns.Item(1).ColumnWidth = 22.44140625

# The newly added ManageNewsPage becomes the active/selected tab.
$newsSheet.Activate()
